$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("prueba" / "Cervezas") - remaining rows shift up automatically.
$ws.Rows.Item(2).Delete()

# Update the "Cantidad" (quantity) column for the remaining products
# to reflect the new (low) stock levels.
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 8
$ws.Range("C5").Value = 10

# "jack daniels" is now low on stock / marked inactive.
$ws.Range("D5").Value = "inactivo"
